# Apply trade #96 closing update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.71
$summary.Range("B4").Value = -0.28
$summary.Range("B5").Value = -0.06
$summary.Range("B6").Value = 96
$summary.Range("B7").Value = 40
$summary.Range("B9").Value = 41.67

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.70999999999999
$status.Range("D4").Value = 96
$status.Range("E4").Value = -0.28
$status.Range("F4").Value = -0.29
$status.Range("G4").Value = 41.67

# ---------------------------------------------------------------------------
# New closed trade row appended to both the "All Trades" and
# "MarketMaking" sheets (row 97). The leading apostrophe on the date
# value keeps it a plain text string (matching the existing rows)
# instead of letting Excel auto-convert it to a date serial number.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Cells.Item(97, 1).Value = 96
    $sheet.Cells.Item(97, 2).Value = "'2026-02-17"
    $sheet.Cells.Item(97, 3).Value = "09:09:31"
    $sheet.Cells.Item(97, 4).Value = "MarketMaking"
    $sheet.Cells.Item(97, 5).Value = "UP"
    $sheet.Cells.Item(97, 6).Value = 0.95
    $sheet.Cells.Item(97, 7).Value = 0.98
    $sheet.Cells.Item(97, 8).Value = "CLOSED"
    $sheet.Cells.Item(97, 9).Value = 3.1579
    $sheet.Cells.Item(97, 10).Value = 0.03
    $sheet.Cells.Item(97, 11).Value = 99.70999999999999
    $sheet.Cells.Item(97, 12).Value = 0
    $sheet.Cells.Item(97, 13).Value = 0
    $sheet.Cells.Item(97, 14).Value = 0.6
    $sheet.Cells.Item(97, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(97, 16).Value = "early_exit"
    $sheet.Cells.Item(97, 17).Value = 0.13
}
